$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-25 Tuesday", "2024-06-26 Wednesday"),
    @("925×9=", "835×7="),
    @("287×4=", "995×9="),
    @("992×5=", "962×6="),
    @("212×7=", "764×2="),
    @("930×3=", "998×9="),
    @("846×9=", "325×6="),
    @("260×3=", "669×9="),
    @("165×3=", "599×6="),
    @("325×4=", "610×9="),
    @("821×4=", "443×4="),
    @("413×2=", "699×7="),
    @("935×3=", "395×9="),
    @("133×6=", "274×4="),
    @("696×3=", "648×9="),
    @("621×2=", "725×3="),
    @("296×6=", "122×4="),
    @("766×4=", "849×9="),
    @("832×5=", "300×3="),
    @("782×6=", "332×4="),
    @("940×7=", "870×8="),
    @("312×8=", "674×2="),
    @("984×9=", "822×7="),
    @("776×2=", "283×7="),
    @("818×9=", "916×6="),
    @("798×4=", "310×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
